# Suivi stage.docx - "Debut sprint 3, creation de compte et modif"
#
# The last paragraph of the document ("Entretien avec le tuteur : ...
# administrateur du site.") currently ends with the hidden _GoBack
# bookmark. The edit appends a full new day's entry (31/05/2018) after
# it, and the _GoBack bookmark moves into the new "Debut " paragraph.

$d = $word.ActiveDocument

# The _GoBack bookmark currently sits at the end of the last paragraph;
# it will be re-created further down in the new content, so drop the
# old one first (Word only ever keeps a single _GoBack).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Start a fresh, empty paragraph right after the current last paragraph
# so the existing content is left completely untouched.
$lastPar = $d.Paragraphs.Last
$tail = $lastPar.Range
$tail.Collapse(2)
$tail.InsertParagraphAfter()

# Replace that (empty) new paragraph's range with the whole block of new
# paragraphs at once, expressed as raw WordprocessingML.
$newPar = $d.Paragraphs.Last
$insertionPoint = $newPar.Range

$newContentXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>31/05/2018</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Début</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>du sprint 3, réalisation/reprise du menu de création de comptes pour la partie back office.</w:t></w:r><w:r><w:t xml:space="preserve"> Mise en place de deux types de comptes : administrateur et normal.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Ce qui a été fait :</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Connexion</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Modification de comptes pour un administrateur (changement de statut ou de MDP)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Création de comptes pour un administrateur</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Mise en place du super administrateur </w:t></w:r><w:r><w:t xml:space="preserve">(administrateur </w:t></w:r><w:r><w:t>impossible à supprimer</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>
'@

$insertionPoint.InsertXML($newContentXml)

Write-Host "Paragraphs count:" $d.Paragraphs.Count
